$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update D23: "Nee" -> "Ja"
$ws.Range("D23").Value = "Ja"

# Update C26: blank -> 0.2
$ws.Range("C26").Value = 0.2

# Update D26: "Nee" -> "Ja"
$ws.Range("D26").Value = "Ja"

# Update selection to D24
$ws.Range("D24").Select()
